$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B, D, E, F, G between each pair of rows that
# represent the same item (same column C) but whose stock figures
# were transposed between the two rows.
$rowPairs = @(
    @(109, 110),
    @(145, 146),
    @(152, 153),
    @(162, 163),
    @(175, 176),
    @(225, 226),
    @(369, 370),
    @(404, 405),
    @(409, 410),
    @(445, 446),
    @(465, 466),
    @(511, 512),
    @(533, 534),
    @(535, 536),
    @(618, 619),
    @(714, 715),
    @(774, 775),
    @(776, 777),
    @(807, 808),
    @(833, 834),
    @(835, 836),
    @(837, 838),
    @(887, 888),
    @(1032, 1033)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $row1 = $pair[0]
    $row2 = $pair[1]
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $val1 = $ws.Range($addr1).Value2
        $val2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value2 = $val2
        $ws.Range($addr2).Value2 = $val1
    }
}

Write-Host "Done swapping $($rowPairs.Count) row pairs."